$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Marking" row (row 11): Right marking and Wrong marking values
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Update "Total" row (row 12): recalculated totals and score string
$ws.Range("B12").Value = 216
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "214/252"
